$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The measured washer thickness / offset (B8) was slightly off; it is now
# derived from the actual micrometer reading (12.7 mm diameter / 2) rather
# than the previously hard-coded value of 6.5.
$ws.Range("B8").Formula = "=12.7 / 2"

# Reflect the cursor/selection position as last left by the editor.
$ws.Range("G18").Select()
